$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet 1: "API Design for LRS"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("API Design for LRS")

# Row 21 / column E: login response JSON gains rewards_earned / rewards_spent /
# rewards_balance fields inside wallet_info (reproducing the original typo).
$loginRs = @'
{
  "loyalty_end_user_login_rs": {
    "user_info": {
      "status": "success",
      "message": "Login successful",
      "user_info": {
        "user_id": "101",
        "email": "user@example.com",
        "tier": {
          "tier_id": "3",
          "tier_name": "Gold"
        },
        "assigned_offers": [
          {
            "offer_id": "1",
            "offer_name": "Exclusive Access",
            "offer_desc": "Exclusive product launch preview"
          },
          {
            "offer_id": "2",
            "offer_name": "Priority Support",
            "offer_desc": "Enjoy priority access to our customer support team"
          }
        ],
        "wallet_info": {
          "ada_balance": "1200",
          "rewards_earned":"800",
         "rewards_spent:"300",
          "rewards_balance":"500",
          "transactions": [
            {
              "transaction_id": "1",
              "date": "2024-01-10",
              "amount": "100",
              "type": "credit",
              "desc": "ADA reward for Gold-tier spending"
            },
            {
              "transaction_id": "2",
              "date": "2024-01-15",
              "amount": "50",
              "type": "debit",
              "desc": "Purchase of product"
            }
          ]
        }
      }
    }
  }
}
'@
$ws1.Range("E21").Value2 = $loginRs

# Row 25: "create a wallet" story becomes "integrate existing wallet".
$ws1.Range("A25").Value2 = "After successful login, the end user needs to integrate existing wallet"

$walletIntegrateRq = @'
{
  "loyalty_end_user_wallet_integrate_rq": {
    "header": {
      "user_name": "endUser",
      "product": "lrs",
      "request_type": "END_USER_CREATE_WALLET"
    },
    "wallet_info": {
      "user_id": "1",
      "currency_type": "ada",
      "wallet_name": "cardanoWallet",
     "wallet_address": "addr_test1vppvktxxw8eyhwkdf1jzq5xxqpxk8sj9d7pzvntfkng94ycn9mjxq",
    }
  }
}
'@
$ws1.Range("D25").Value2 = $walletIntegrateRq

$walletIntegrateRs = @'
{
  "loyalty_end_user_wallet_integrate_rs": {
     "status": "success"
  }
}
'@
$ws1.Range("E25").Value2 = $walletIntegrateRs

# The old "data is shown once" comment cell is removed entirely.
$ws1.Range("F25").ClearContents()

# Row grew taller to fit the new response text.
$ws1.Rows.Item(25).RowHeight = 245.25

# Selection ends up parked on the edited cell.
$ws1.Activate()
$ws1.Range("E21").Select()

# ----------------------------------------------------------------------
# Sheet 2: "Integration with HBS"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Integration with HBS")

$ws2.Range("A3").Value2 = "Insert Hotel master data in lrs"
$ws2.Range("A5").Value2 = "Once a guest completes a payment (or a set of guests at configurable schedule as per HBS) in the HBS, their booking information (name, email, total payment, currency) is sent to the lrs"

$ws2.Activate()
$ws2.Range("D4").Select()

# Restore the first sheet as the active one, matching the saved view state.
$ws1.Activate()
